# The post corresponding to row 765 ("「ひらめき」") was removed from the
# underlying data source, so its entire row is deleted from the sheet.
# Deleting the row shifts every following row up by one, which also matches
# the updated sheet dimension (A1:C873 -> A1:C872).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(765).Delete()
